$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that need to remain plain text even though they look numeric/percent
# (prices in column D, percentages in column E). Force Text number format,
# set literal text value, then reset style to Normal so no stray style sticks.

$textCells = @{
    D2 = "244.79"
    E2 = "-0.85%"
    E3 = "2.31%"
    D4 = "5.074"
    E4 = "-0.11%"
    D5 = "0.05685"
    E5 = "1.25%"
    E6 = "-0.40%"
    D7 = "0.8208"
    E7 = "0.96%"
    D8 = "0.8384"
    E8 = "-0.70%"
    D9 = "0.1328"
    E9 = "-1.32%"
    D10 = "0.06908"
    E10 = "-0.95%"
    D11 = "0.02849"
    E11 = "-0.95%"
    D12 = "0.09398"
    E12 = "-0.18%"
    D13 = "0.001531"
    E13 = "-0.04%"
    D14 = "0.04097"
    E14 = "-12.21%"
    D15 = "0.0005998"
    E15 = "0.14%"
    D16 = "0.006123"
    E16 = "-1.54%"
    D17 = "3.510"
    E17 = "-2.34%"
    D18 = "3.002"
    E18 = "-0.30%"
    D19 = "2.310"
    E19 = "9.05%"
    E20 = "-0.21%"
    D21 = "0.03178"
    E21 = "-0.21%"
    D22 = "0.1296"
    E22 = "-1.80%"
    D23 = "3.573"
    D24 = "0.1374"
    E24 = "1.78%"
    D25 = "0.001217"
    E25 = "-2.81%"
    D26 = "0.003952"
    E26 = "-14.10%"
    D27 = "0.00009795"
    E27 = "2.02%"
    D28 = "0.0001938"
    E28 = "-0.02%"
    D40 = "0.03690"
    E40 = "0.28%"
    D41 = "0.005837"
    E41 = "-5.86%"
    D42 = "0.1053"
    E42 = "-21.87%"
    D43 = "0.002342"
    E43 = "-6.33%"
    D44 = "0.009374"
    E44 = "5.30%"
    D45 = "0.00005195"
    E45 = "-1.92%"
    D47 = "0.1015"
    E47 = "-32.30%"
    D48 = "0.002591"
    E48 = "2.57%"
}

foreach ($ref in $textCells.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $textCells[$ref]
    $cell.Style = "Normal"
}

# Plain text cells (coin names / links) - safe to assign directly
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
